$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") is stored as text in this sheet (values use "." as a thousands
# separator, e.g. "47.740.35", so it can never be a real number). Force text format on
# every Price cell we touch so Excel does not reinterpret the typed value as a number
# (which would silently drop things like trailing zeros: "316.00" -> 316).
$priceCells = @("D2","D3","D5","D6","D7","D8","D10","D11","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D29","D30","D31","D32","D34","D35","D36","D38","D39","D40","D42","D43","D45","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Cell value updates (from the Feb 12 2024 11:39:20 UTC GitHub Actions refresh)
$ws.Range("D2").Value = "47.740.35"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "2.478.26"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "316.00"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").Value = "103.72"
$ws.Range("E6").Value = "  -5.36%  "
$ws.Range("D7").Value = "0.516"
$ws.Range("E7").Value = "  -3.14%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -3.95%  "
$ws.Range("D10").Value = "38.49"
$ws.Range("E10").Value = "  -4.95%  "
$ws.Range("D11").Value = "20.44"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "7.00"
$ws.Range("E14").Value = "  -3.90%  "
$ws.Range("D15").Value = "2.865.79"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("D16").Value = "2.483.59"
$ws.Range("E16").Value = "  -1.58%  "
$ws.Range("D17").Value = "0.820"
$ws.Range("E17").Value = "  -4.12%  "
$ws.Range("D18").Value = "47.668.65"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "2.91"
$ws.Range("E19").Value = "  +7.71%  "
$ws.Range("D20").Value = "12.63"
$ws.Range("E20").Value = "  -6.25%  "
$ws.Range("D21").Value = "6.50"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("D22").Value = "0.0₃0924"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").Value = "278.23"
$ws.Range("E23").Value = "  +5.09%  "
$ws.Range("D24").Value = "70.63"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  -3.66%  "
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").Value = "9.52"
$ws.Range("E29").Value = "  -5.77%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.135"
$ws.Range("E30").Value = "  -5.50%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "34.17"
$ws.Range("E31").Value = "  -4.63%  "
$ws.Range("D32").Value = "49.13"
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").Value = "18.84"
$ws.Range("E34").Value = "  -4.65%  "
$ws.Range("D35").Value = "5.23"
$ws.Range("E35").Value = "  -3.32%  "
$ws.Range("D36").Value = "0.0764"
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("E37").Value = "  -3.60%  "
$ws.Range("D38").Value = "4.46"
$ws.Range("E38").Value = "  -5.64%  "
$ws.Range("D39").Value = "2.83"
$ws.Range("E39").Value = "  -5.84%  "
$ws.Range("D40").Value = "122.89"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("E41").Value = "  -1.86%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "21.64"
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").Value = "1.987.95"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("E46").Value = "  -2.02%  "
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("D49").Value = "8.86"
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("D50").Value = "5.04"
$ws.Range("E50").Value = "  -3.54%  "
$ws.Range("D51").Value = "78.84"
$ws.Range("E51").Value = "  -0.53%  "
